# Revert "Powerpoint writer: consolidate text run nodes."
#
# The title text runs "First " (slide 1) and "Third " (slide 3) were
# previously merged with the following "slide" run's leading space into
# a single trailing-space run. This splits the trailing space back out
# into its own run, so each heading reads as three runs:
#   "First"/" "/"slide"   and   "Third"/" "/"slide"
#
# Re-assigning .Text on a Characters() sub-range that only covers part
# of an existing run forces PowerPoint to split that run, which is
# exactly the OOXML shape the diff calls for.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 5).Text = "First"

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 5).Text = "Third"
